$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The observation records in rows 6-10 were re-ordered (a single 5-way
# rotation): the new row 6 gets the old row 9's data, new row 7 gets old
# row 8's data, new row 8 gets old row 6's data, new row 9 gets old row
# 10's data and new row 10 gets old row 7's data. Capture every source
# row's data into variables BEFORE writing anything back, since the
# rotation is one big cycle touching all five rows.

$data6  = $ws.Range("A6:H6").Value2
$data7  = $ws.Range("A7:H7").Value2
$data8  = $ws.Range("A8:H8").Value2
$data9  = $ws.Range("A9:H9").Value2
$data10 = $ws.Range("A10:H10").Value2

$q6  = $ws.Range("Q6").Value2
$r6  = $ws.Range("R6").Value2
$q7  = $ws.Range("Q7").Value2
$r7  = $ws.Range("R7").Value2
$q8  = $ws.Range("Q8").Value2
$r8  = $ws.Range("R8").Value2
$q9  = $ws.Range("Q9").Value2
$r9  = $ws.Range("R9").Value2
$q10 = $ws.Range("Q10").Value2
$r10 = $ws.Range("R10").Value2

# Write back the rotated data (columns A-H).
$ws.Range("A6:H6").Value   = $data9
$ws.Range("A7:H7").Value   = $data8
$ws.Range("A8:H8").Value   = $data6
$ws.Range("A9:H9").Value   = $data10
$ws.Range("A10:H10").Value = $data7

# The coordinate columns (Q, R) rotate the same way, but the new values
# are also rounded to whole numbers.
$ws.Range("Q6").Value  = [math]::Round($q9, 0)
$ws.Range("R6").Value  = [math]::Round($r9, 0)
$ws.Range("Q7").Value  = [math]::Round($q8, 0)
$ws.Range("R7").Value  = [math]::Round($r8, 0)
$ws.Range("Q8").Value  = [math]::Round($q6, 0)
$ws.Range("R8").Value  = [math]::Round($r6, 0)
$ws.Range("Q9").Value  = [math]::Round($q10, 0)
$ws.Range("R9").Value  = [math]::Round($r10, 0)
$ws.Range("Q10").Value = [math]::Round($q7, 0)
$ws.Range("R10").Value = [math]::Round($r7, 0)

# Column L (empty marker cell) moves from row 6 to row 8. Copy it (rather
# than just writing "") so the destination keeps the same empty-text cell
# shape the source had, instead of becoming a fully blank cell.
$ws.Range("L6").Copy($ws.Range("L8"))
$ws.Range("L6").ClearContents()

# The "Starttid"/"Sluttid" time-of-day columns (Z, AB) - both always
# "00:00" - are removed entirely for every one of these rows.
$ws.Range("Z6:Z10").ClearContents()
$ws.Range("AB6:AB10").ClearContents()
